# Generate Report for Handback
#
# Marks both localization target sheets (zh-cn, de-de) as handed back:
#   - Status column (C) -> "Handed back: in sync with en-US"
#   - New "Latest Target File" (F) / "Latest Handback File" (G) hyperlink
#     cells mirroring the Source File Name (A) / Latest Handoff File (D)
#     entries for each row.
#   - "Latest Handback DateTime" (H) stamped with the handback timestamp
#     (differs per-language since each language finished at a different time).

$wb = $excel.ActiveWorkbook

function Set-HandbackRow($ws, $row, $mdName, $xlfName, $mdUrl, $xlfUrl, $handbackDateTime) {
    # Status -> handed back, in sync with the English source.
    $ws.Cells.Item($row, 3).Value = "Handed back: in sync with en-US"

    # F = "Latest Target File" -> same markdown file as column A, now
    # also linked from the handback side.
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $mdUrl, "", "", $mdName)

    # G = "Latest Handback File" -> same xlf file as column D, linked
    # from the handback side.
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $xlfUrl, "", "", $xlfName)

    # H = "Latest Handback DateTime"
    $ws.Cells.Item($row, 8).Value = $handbackDateTime
}

$mdName1 = "6f3372ea-768d-46e8-a856-8c846f7c8ac0.md"
$mdName2 = "98c37869-d990-467a-86a4-020f4d10662c.md"
$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/281b5a87328ccf46942c6b50d193d37b6c5efcd0/e2e/6f3372ea-768d-46e8-a856-8c846f7c8ac0.md"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/281b5a87328ccf46942c6b50d193d37b6c5efcd0/e2e/98c37869-d990-467a-86a4-020f4d10662c.md"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$xlfZh1 = "6f3372ea-768d-46e8-a856-8c846f7c8ac0.5bd438904a775b4547bf051991ae2f98a0d454fb.zh-cn.xlf"
$xlfZh2 = "98c37869-d990-467a-86a4-020f4d10662c.f41a6e0425d6d4e81b5b8fc2a66cf4ced1955376.zh-cn.xlf"
$xlfZhUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6d09383f7c2779e82f3c782f0a48700e7a6da849/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/high/$xlfZh1"
$xlfZhUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6d09383f7c2779e82f3c782f0a48700e7a6da849/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/high/$xlfZh2"

Set-HandbackRow $wsZh 2 $mdName1 $xlfZh1 $mdUrl1 $xlfZhUrl1 "2016-03-20 05:37:49"
Set-HandbackRow $wsZh 3 $mdName2 $xlfZh2 $mdUrl2 $xlfZhUrl2 "2016-03-20 05:37:49"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$xlfDe1 = "6f3372ea-768d-46e8-a856-8c846f7c8ac0.5bd438904a775b4547bf051991ae2f98a0d454fb.de-de.xlf"
$xlfDe2 = "98c37869-d990-467a-86a4-020f4d10662c.f41a6e0425d6d4e81b5b8fc2a66cf4ced1955376.de-de.xlf"
$xlfDeUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fd299a1570a1b7183b04b8d891095ebc1c4a7387/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/high/$xlfDe1"
$xlfDeUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fd299a1570a1b7183b04b8d891095ebc1c4a7387/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/high/$xlfDe2"

Set-HandbackRow $wsDe 2 $mdName1 $xlfDe1 $mdUrl1 $xlfDeUrl1 "2016-03-20 05:38:02"
Set-HandbackRow $wsDe 3 $mdName2 $xlfDe2 $mdUrl2 $xlfDeUrl2 "2016-03-20 05:38:02"
